$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$used = $ws.UsedRange
$rows = $used.Rows.Count
Write-Host "Rows: $rows"
for ($r = 2; $r -le $rows; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2
    if ($val -ne $null -and $val.ToString().Contains(",")) {
        $parts = $val.ToString().Split(",")
        $n = $parts.Length
        for ($i = 0; $i -lt $n; $i++) {
            $parts[$i] = $parts[$i].Trim()
        }
        $revParts = $parts[($n - 1)..0]
        $newVal = [String]::Join(", ", $revParts)
        $cell.Value2 = $newVal
    }
}
